$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 530 (existing rows 530-564 shift down to 533-567),
# matching the weekly refresh that prepends a new "Especial/Primera/Segunda"
# price-report group for the latest date (44585) ahead of the prior groups.
$ws.Rows("530:532").Insert()

# Row 530 - Especial
$ws.Range("A530").Value = 8
$ws.Range("B530").Value = "Terminal La Palmera de La Serena"
$ws.Range("C530").Value = "Coquimbo"
$ws.Range("D530").Value = 44585
$ws.Range("E530").Value = 4
$ws.Range("F530").Value = "Fruta"
$ws.Range("G530").Value = 100101
$ws.Range("H530").Value = "Berries"
$ws.Range("I530").Value = 100112025
$ws.Range("J530").Value = "Frutilla"
$ws.Range("K530").Value = "Sin especificar"
$ws.Range("L530").Value = "Especial"
$ws.Range("M530").Value = 500
$ws.Range("N530").Value = 11500
$ws.Range("O530").Value = 12000
$ws.Range("P530").Value = 11750
$ws.Range("Q530").Value = "`$/bandeja 7 kilos"
$ws.Range("R530").Value = "Provincia de Melipilla"
$ws.Range("S530").Value = 1679
$ws.Range("T530").Value = 7

# Row 531 - Primera
$ws.Range("A531").Value = 8
$ws.Range("B531").Value = "Terminal La Palmera de La Serena"
$ws.Range("C531").Value = "Coquimbo"
$ws.Range("D531").Value = 44585
$ws.Range("E531").Value = 4
$ws.Range("F531").Value = "Fruta"
$ws.Range("G531").Value = 100101
$ws.Range("H531").Value = "Berries"
$ws.Range("I531").Value = 100112025
$ws.Range("J531").Value = "Frutilla"
$ws.Range("K531").Value = "Sin especificar"
$ws.Range("L531").Value = "Primera"
$ws.Range("M531").Value = 400
$ws.Range("N531").Value = 9500
$ws.Range("O531").Value = 10000
$ws.Range("P531").Value = 9750
$ws.Range("Q531").Value = "`$/bandeja 7 kilos"
$ws.Range("R531").Value = "Provincia de Melipilla"
$ws.Range("S531").Value = 1393
$ws.Range("T531").Value = 7

# Row 532 - Segunda
$ws.Range("A532").Value = 8
$ws.Range("B532").Value = "Terminal La Palmera de La Serena"
$ws.Range("C532").Value = "Coquimbo"
$ws.Range("D532").Value = 44585
$ws.Range("E532").Value = 4
$ws.Range("F532").Value = "Fruta"
$ws.Range("G532").Value = 100101
$ws.Range("H532").Value = "Berries"
$ws.Range("I532").Value = 100112025
$ws.Range("J532").Value = "Frutilla"
$ws.Range("K532").Value = "Sin especificar"
$ws.Range("L532").Value = "Segunda"
$ws.Range("M532").Value = 400
$ws.Range("N532").Value = 7500
$ws.Range("O532").Value = 8000
$ws.Range("P532").Value = 7750
$ws.Range("Q532").Value = "`$/bandeja 7 kilos"
$ws.Range("R532").Value = "Provincia de Melipilla"
$ws.Range("S532").Value = 1107
$ws.Range("T532").Value = 7
